$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last_edited_time" column (D) values to reflect the latest
# Notion sync timestamps for the "khach hang" report rows.
$ws.Range("D2:D15").Value = "2024-07-04T09:36:00.000Z"
$ws.Range("D16:D24").Value = "2024-07-04T09:33:00.000Z"
$ws.Range("D25:D67").Value = "2024-07-04T09:34:00.000Z"
$ws.Range("D68:D72").Value = "2024-07-04T09:35:00.000Z"
